$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are plain decimal numbers (e.g. "9.00", "573.00").
# Assigning them directly would make Excel auto-detect them as numeric values and
# drop the significant trailing zeros / leading-zero formatting used throughout the
# "Price" column, which is stored as literal text. Mark those specific cells as Text
# before writing so the literal string is preserved exactly, then drop back to the
# default "Normal" style afterwards so no stray formatting is left behind.
$textCells = @("D5", "D6", "D8", "D10", "D12", "D13", "D15", "D16", "D20", "D21", "D23", "D25", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D36", "D41", "D42", "D43", "D44", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "71.677.90"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "3.687.32"
$ws.Range("E3").Value = "  +8.49%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "589.17"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "179.75"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "3.677.98"
$ws.Range("E7").Value = "  +8.42%  "
$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  +4.73%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.202"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +4.54%  "
$ws.Range("D12").Value = "49.95"
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("D13").Value = "0.0000287"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "4.280.53"
$ws.Range("E14").Value = "  +8.44%  "
$ws.Range("D15").Value = "685.11"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "9.00"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "71.683.62"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.672.84"
$ws.Range("E18").Value = "  +8.08%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "18.10"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "11.66"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("D23").Value = "6.31"
$ws.Range("E23").Value = "  +17.42%  "
$ws.Range("E24").Value = "  +4.23%  "
$ws.Range("D25").Value = "103.96"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("E26").Value = "  +3.79%  "
$ws.Range("D27").Value = "2.85"
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("D28").Value = "10.20"
$ws.Range("D29").Value = "35.47"
$ws.Range("E29").Value = "  +5.89%  "
$ws.Range("D30").Value = "9.23"
$ws.Range("E30").Value = "  +5.35%  "
$ws.Range("E31").Value = "  +6.34%  "
$ws.Range("D32").Value = "4.26"
$ws.Range("E32").Value = "  +12.69%  "
$ws.Range("D33").Value = "573.00"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Value = "11.31"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("D36").Value = "59.46"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("D37").Value = "3.778.97"
$ws.Range("E37").Value = "  +4.72%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").Value = "0.0₃0778"
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("D41").Value = "35.44"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  +5.66%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0462"
$ws.Range("E43").Value = "  +8.34%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.79"
$ws.Range("E44").Value = "  +3.35%  "
$ws.Range("E45").Value = "  +4.87%  "
$ws.Range("E46").Value = "  +7.64%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "134.59"
$ws.Range("E51").Value = "  +2.72%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
